{"js": "// Replace the two-digit multiplication problems throughout the document.\n// Each \"AxB=\" expression is unique in the document, so we can safely\n// search for the exact old text and replace it with the new text while\n// keeping the existing run formatting intact.\nconst replacements = [\n  [\"20\u00d779=\", \"70\u00d766=\"],\n  [\"63\u00d733=\", \"58\u00d728=\"],\n  [\"38\u00d799=\", \"32\u00d723=\"],\n  [\"73\u00d711=\", \"22\u00d745=\"],\n  [\"39\u00d734=\", \"76\u00d722=\"],\n  [\"37\u00d741=\", \"74\u00d796=\"],\n  [\"52\u00d746=\", \"40\u00d777=\"],\n  [\"67\u00d727=\", \"92\u00d755=\"],\n  [\"16\u00d765=\", \"49\u00d722=\"],\n  [\"14\u00d767=\", \"90\u00d769=\"],\n  [\"18\u00d749=\", \"39\u00d745=\"],\n  [\"56\u00d746=\", \"65\u00d793=\"],\n  [\"92\u00d741=\", \"60\u00d731=\"],\n  [\"72\u00d753=\", \"56\u00d737=\"],\n  [\"20\u00d750=\", \"54\u00d744=\"],\n  [\"25\u00d733=\", \"45\u00d748=\"],\n  [\"67\u00d750=\", \"12\u00d751=\"],\n  [\"62\u00d786=\", \"60\u00d723=\"],\n  [\"25\u00d751=\", \"37\u00d759=\"],\n  [\"84\u00d777=\", \"15\u00d755=\"],\n  [\"66\u00d731=\", \"92\u00d752=\"],\n  [\"93\u00d711=\", \"46\u00d794=\"],\n  [\"28\u00d760=\", \"57\u00d723=\"],\n  [\"44\u00d713=\", \"42\u00d734=\"],\n  [\"48\u00d719=\", \"90\u00d717=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems throughout the document.\n# Each \"AxB=\" expression is unique in the document, so Find/Replace on the\n# exact old text safely targets only the intended run, preserving all\n# existing character/paragraph formatting.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"20\u00d779=\", \"70\u00d766=\"),\n  @(\"63\u00d733=\", \"58\u00d728=\"),\n  @(\"38\u00d799=\", \"32\u00d723=\"),\n  @(\"73\u00d711=\", \"22\u00d745=\"),\n  @(\"39\u00d734=\", \"76\u00d722=\"),\n  @(\"37\u00d741=\", \"74\u00d796=\"),\n  @(\"52\u00d746=\", \"40\u00d777=\"),\n  @(\"67\u00d727=\", \"92\u00d755=\"),\n  @(\"16\u00d765=\", \"49\u00d722=\"),\n  @(\"14\u00d767=\", \"90\u00d769=\"),\n  @(\"18\u00d749=\", \"39\u00d745=\"),\n  @(\"56\u00d746=\", \"65\u00d793=\"),\n  @(\"92\u00d741=\", \"60\u00d731=\"),\n  @(\"72\u00d753=\", \"56\u00d737=\"),\n  @(\"20\u00d750=\", \"54\u00d744=\"),\n  @(\"25\u00d733=\", \"45\u00d748=\"),\n  @(\"67\u00d750=\", \"12\u00d751=\"),\n  @(\"62\u00d786=\", \"60\u00d723=\"),\n  @(\"25\u00d751=\", \"37\u00d759=\"),\n  @(\"84\u00d777=\", \"15\u00d755=\"),\n  @(\"66\u00d731=\", \"92\u00d752=\"),\n  @(\"93\u00d711=\", \"46\u00d794=\"),\n  @(\"28\u00d760=\", \"57\u00d723=\"),\n  @(\"44\u00d713=\", \"42\u00d734=\"),\n  @(\"48\u00d719=\", \"90\u00d717=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Forward = $true\n  $find.Wrap = 1\n\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
